# Auto-generated: apply scheduled-runner market data refresh to Golem_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2101.5833
$ws.Range("I40").Value = 1965.5264
$ws.Range("J40").Value = 2618.6
$ws.Range("K40").Value = 1965.5264
$ws.Range("L40").Value = 2618.6
$ws.Range("M40").Value = -1790.5264
$ws.Range("N40").Value = -2968.6
$ws.Range("H51").Value = 115399.8
$ws.Range("J51").Value = 115399.8
$ws.Range("L51").Value = 115399.8
$ws.Range("N51").Value = -116367.8
$ws.Range("H92").Value = 525.25
$ws.Range("I92").Value = 525.25
$ws.Range("K92").Value = 525.25
$ws.Range("M92").Value = 722.75
$ws.Range("H132").Value = 2125.4285
$ws.Range("I132").Value = 1520.6
$ws.Range("J132").Value = 3637.5
$ws.Range("K132").Value = 4561.799999999999
$ws.Range("L132").Value = 10912.5
$ws.Range("M132").Value = -2031.799999999999
$ws.Range("N132").Value = -15972.5
$ws.Range("H138").Value = 3684.625
$ws.Range("I138").Value = 3872.25
$ws.Range("J138").Value = 3497
$ws.Range("K138").Value = 11616.75
$ws.Range("L138").Value = 10491
$ws.Range("M138").Value = -6476.75
$ws.Range("N138").Value = -20771
$ws.Range("H141").Value = 3438
$ws.Range("I141").Value = 2797.5
$ws.Range("K141").Value = 8392.5
$ws.Range("M141").Value = -3212.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").Value = ""
$ws.Range("H132").Value = 1977.1428
$ws.Range("I132").Value = 1977.1428
$ws.Range("K132").Value = 5931.428400000001
$ws.Range("M132").Value = -3401.428400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2383.7693
$ws.Range("I94").Value = 1998.625
$ws.Range("K94").Value = 1998.625
$ws.Range("M94").Value = -1547.625
$ws.Range("H107").Value = 1470.25
$ws.Range("I107").Value = 1466
$ws.Range("J107").Value = 1500
$ws.Range("K107").Value = 1466
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = 454
$ws.Range("N107").Value = -5340
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").Value = ""
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").Value = ""
$ws.Range("H134").Value = 1567.3334
$ws.Range("I134").Value = 1567.3334
$ws.Range("K134").Value = 4702.0002
$ws.Range("M134").Value = -2167.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1306.5918
$ws.Range("I4").Value = 621.2646999999999
$ws.Range("K4").Value = 1863.7941
$ws.Range("M4").Value = -1751.7941
$ws.Range("H129").Value = 1120
$ws.Range("I129").Value = 493.33334
$ws.Range("J129").Value = 3000
$ws.Range("K129").Value = 1480.00002
$ws.Range("L129").Value = 9000
$ws.Range("M129").Value = 3519.99998
$ws.Range("N129").Value = -19000
$ws.Range("H131").Value = 2360.1428
$ws.Range("I131").Value = 1445.4445
$ws.Range("J131").Value = 4006.6
$ws.Range("K131").Value = 4336.333500000001
$ws.Range("L131").Value = 12019.8
$ws.Range("M131").Value = 703.6664999999994
$ws.Range("N131").Value = -22099.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 3685882.8
$ws.Range("J11").Value = 232221.78
$ws.Range("L11").Value = 232221.78
$ws.Range("N11").Value = -232499.78
$ws.Range("H97").Value = 1853.4546
$ws.Range("I97").Value = 1509.7778
$ws.Range("J97").Value = 3400
$ws.Range("K97").Value = 1509.7778
$ws.Range("L97").Value = 3400
$ws.Range("M97").Value = -1013.7778
$ws.Range("N97").Value = -4392
$ws.Range("H102").Value = 51099.375
$ws.Range("I102").Value = 58113.57
$ws.Range("K102").Value = 58113.57
$ws.Range("M102").Value = -56491.57

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 25713.637
$ws.Range("I40").Value = 15871.667
$ws.Range("K40").Value = 15871.667
$ws.Range("M40").Value = -15735.667
$ws.Range("H44").Value = 799999
$ws.Range("J44").Value = 799999
$ws.Range("L44").Value = 799999
$ws.Range("N44").Value = -800911
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").Value = ""
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").Value = ""
$ws.Range("H93").Value = 3216.6667
$ws.Range("I93").Value = 3000
$ws.Range("J93").Value = 3325
$ws.Range("K93").Value = 3000
$ws.Range("L93").Value = 3325
$ws.Range("M93").Value = -1752
$ws.Range("N93").Value = -5821
$ws.Range("H97").Value = 18000
$ws.Range("J97").Value = 18000
$ws.Range("L97").Value = 18000
$ws.Range("N97").Value = -19982
$ws.Range("H132").Value = 994.6667
$ws.Range("I132").Value = 994.6667
$ws.Range("K132").Value = 2984.0001
$ws.Range("M132").Value = -454.0001000000002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 12000
$ws.Range("I45").Value = 12000
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 12000
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -11509
$ws.Range("N45").Value = ""
$ws.Range("H100").Value = 686.125
$ws.Range("I100").Value = 548.3333
$ws.Range("J100").Value = 1099.5
$ws.Range("K100").Value = 1096.6666
$ws.Range("L100").Value = 2199
$ws.Range("M100").Value = -555.6666
$ws.Range("N100").Value = -3281
$ws.Range("H132").Value = 2428.4285
$ws.Range("I132").Value = 1899.8
$ws.Range("J132").Value = 3750
$ws.Range("K132").Value = 5699.4
$ws.Range("L132").Value = 11250
$ws.Range("M132").Value = -3169.4
$ws.Range("N132").Value = -16310
$ws.Range("H136").Value = 11582.333
$ws.Range("I136").Value = 12069.857
$ws.Range("K136").Value = 36209.571
$ws.Range("M136").Value = -33659.571
$ws.Range("H141").Value = 90000
$ws.Range("J141").Value = 90000
$ws.Range("L141").Value = 90000
$ws.Range("N141").Value = -100360
